# "updated directors and ec"
# Update the names of the Executive Committee / Directors in column A.
# (Titles in column B are unchanged.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value  = "Tommy Wunsch"    # was Jacob Lee   (VP of Finance)
$ws.Range("A6").Value  = "Elaine Nguyen"   # was Ishan Malik (VP of Chapter Operations)
$ws.Range("A7").Value  = "Brennan Kim"     # was Jason Henkel (Chancellor)
$ws.Range("A8").Value  = "Lohit Potnuru"   # was Nithin Senthil (VP of Professional Activities)
$ws.Range("A9").Value  = "Jessica Lin"     # was Juliana Lee (VP of Community Service)
$ws.Range("A10").Value = "Melanie Sagun"   # was a stray duplicate "Mirsab Mirza" (VP of Scholarship and Awards)
$ws.Range("A11").Value = "Isaac Martinez"  # was Diana Huynh (VP of Alumni Relations)

# Column A/B default alignment normalized to left (matches rest of the sheet).
$ws.Columns("A:B").HorizontalAlignment = -4131
